$wb = $excel.ActiveWorkbook

$wsBAU = $wb.Worksheets.Item("RQSD-BRQSD")
$wsBAU.Range("B4").Value = 1

$wsRQSD = $wb.Worksheets.Item("RQSD-RQSD")
$wsRQSD.Range("B4").Value = 0
